$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the dSF column (F) values - repulled data / recalculated means
$ws.Range("F6").Value = -6
$ws.Range("F8").Value = -1
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 6
$ws.Range("F13").Value = -5
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -5
